$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value is a plain decimal number (e.g. "228.14") must be
# kept as text (matching the source data, which stores prices as literal strings,
# e.g. "34.028.34" using dots as thousands separators). Mark them as Text before
# assigning so Excel does not auto-convert them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "34.028.34"
$ws.Range("E2").Value = "  +10.64%  "
$ws.Range("D3").Value = "1.813.02"
$ws.Range("E3").Value = "  +7.44%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "228.14"
$ws.Range("E5").Value = "  +3.24%  "
$ws.Range("D6").Value = "0.541"
$ws.Range("E6").Value = "  +3.72%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "30.89"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "47.93"
$ws.Range("E9").Value = "  +8.17%  "
$ws.Range("D10").Value = "0.279"
$ws.Range("E10").Value = "  +5.09%  "
$ws.Range("D11").Value = "0.0667"
$ws.Range("E11").Value = "  +6.71%  "
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").Value = "2.072.26"
$ws.Range("E13").Value = "  +7.21%  "
$ws.Range("D14").Value = "1.791.06"
$ws.Range("E14").Value = "  +5.07%  "
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "34.012.23"
$ws.Range("E16").Value = "  +10.51%  "
$ws.Range("E17").Value = "  -3.87%  "
$ws.Range("D18").Value = "4.27"
$ws.Range("E18").Value = "  +6.85%  "
$ws.Range("D19").Value = "69.13"
$ws.Range("E19").Value = "  +3.85%  "
$ws.Range("D20").Value = "255.82"
$ws.Range("E20").Value = "  +3.34%  "
$ws.Range("D21").Value = "0.0₃0742"
$ws.Range("E21").Value = "  +3.77%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "10.38"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("D24").Value = "4.31"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").Value = "158.97"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").Value = "16.55"
$ws.Range("E27").Value = "  +4.21%  "
$ws.Range("E28").Value = "  +3.28%  "
$ws.Range("D29").Value = "7.05"
$ws.Range("E29").Value = "  +4.96%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "3.81"
$ws.Range("E31").Value = "  +9.12%  "
$ws.Range("D32").Value = "0.0509"
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("E33").Value = "  +5.56%  "
$ws.Range("E34").Value = "  +6.64%  "
$ws.Range("D35").Value = "1.549.10"
$ws.Range("E35").Value = "  +2.34%  "
$ws.Range("E36").Value = "  +3.86%  "
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("D38").Value = "84.07"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "0.0187"
$ws.Range("E39").Value = "  +3.90%  "
$ws.Range("D40").Value = "0.619"
$ws.Range("E40").Value = "  +5.73%  "
$ws.Range("E41").Value = "  +3.35%  "
$ws.Range("D43").Value = "0.902"
$ws.Range("E43").Value = "  +6.34%  "
$ws.Range("E44").Value = "  +4.92%  "
$ws.Range("D45").Value = "0.0524"
$ws.Range("E45").Value = "  +4.70%  "
$ws.Range("E46").Value = "  +3.79%  "
$ws.Range("D47").Value = "1.968.46"
$ws.Range("E47").Value = "  +7.60%  "
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("D50").Value = "52.22"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0122"
$ws.Range("E51").Value = "  +5.77%  "
